$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id=1)
$ws.Range("B2").Value = 'What was the ligh t a symphonic atmosphere?'
$ws.Range("C2").Value = 'The sun dipped below the horizon, casting a golden glow across the vast expanse of the  landscape. The evening breeze carried the fresh scent of pine trees, and the distant sound of a  stream could be heard, gently flowing through the valley. The fading ligh t created a magical  atmosphere, with the sky transitioning through various shades of purple, orange, and pink. It felt  as if time itself had slowed down, allowing for a moment of reflection and serenity.'
$ws.Range("D2").Value = 'the sky transitioning through various shades of purple, orange, and pink'
$ws.Range("E2").Value = 'a symphonic atmosphere'
$ws.Range("F2").Value = 'the sky transitioning through various shades of purple, orange, and pink'

# Row 3 (id=2)
$ws.Range("B3").Value = 'What was the name of the ancient oak tree?'
$ws.Range("C3").Value = 'As darkness  began to fall, the first stars appeared,  twinkling faintly in the sky, and the night slowly took hold  of the world, inviting a deep sense of peace and calm. The quiet was only broken by the  occasional rustle of leaves, as if nature itself was taking a deep breath in the cool air. Deep within the  heart of the forest, there stood an ancient oak tree, towering over the other  vegetation like a silent guardian of the woods.'
$ws.Range("D3").Value = 'an ancient oak tree'
$ws.Range("E3").Value = 'styrofoam'
$ws.Range("F3").Value = 'an ancient oak tree'
$ws.Range("G3").Value = $false

# Row 4 (id=3)
$ws.Range("B4").Value = 'What is the name of the tree that was known to the locals as the "Whispering Giant"?'
$ws.Range("C4").Value = 'It was known to the locals as the "Whispering  Giant," a name that carried with it a sense of mystery and reverence. According to local legend,  the tree had been there for centuries, its roots intertwined with the very earth, its branches  stretching high into the sky. Some believed the tree had mystical powers, able to communicate  with those who listened carefully.'
$ws.Range("D4").Value = 'Whispering  Giant'
$ws.Range("E4").Value = 'styrofoam'
$ws.Range("F4").Value = 'Whispering  Giant'
$ws.Range("G4").Value = $false

# Row 5 (id=4)
$ws.Range("B5").Value = 'What did the whispers appear to be?'
$ws.Range("C5").Value = 'Travelers and seek ers of wisdom would often come from far  and wide, hoping to hear the whispers that were said to reveal forgotten truths or offer guidance  in times of need. Others claimed the oak held secrets of the past, tales of long -lost civilizations  buried beneath its  roots, waiting for the right person to uncover them. As the wind moved  through its leaves, the whispers seemed to come alive, echoing in the minds of those who dared  to listen.'
$ws.Range("D5").Value = 'alive'
$ws.Range("E5").Value = 'a snarky snark'
$ws.Range("F5").Value = 'alive'
$ws.Range("G5").Value = $false

# Row 6 (id=5)
$ws.Range("B6").Value = 'What was the name of the city that was a constant whirl of activity?'
$ws.Range("C6").Value = 'In stark contrast to the quiet forest, the bustling city streets were a consta nt whirl of activity. The  modern skyline of towering skyscrapers loomed above, their glass windows reflecting the bright  lights of the city below. Amidst the towering structures, narrow alleyways and older brick  buildings preserved the charm of a bygone er a. It was a place where the old and new coexisted in  a harmonious yet chaotic dance.'
$ws.Range("D6").Value = 'city streets'
$ws.Range("E6").Value = 'san francisco'
$ws.Range("F6").Value = 'city streets'

# Row 7 (id=6)
$ws.Range("B7").Value = 'What was the smell of street food?'
$ws.Range("C7").Value = 'The streets were always crowded with people, each person  moving with purpose, whether they were heading to work, meeting friends, or simply passing  through. The air was fi lled with the aroma of street food —spicy noodles, sizzling skewers, and  freshly baked pastries —drawing crowds to the vibrant market stalls that lined the sidewalks. Everywhere you looked, there was color and life, from the brightly painted storefronts to t he  lively chatter of vendors and customers.'
$ws.Range("D7").Value = 'fi lled'
$ws.Range("E7").Value = 'a cigarette'
$ws.Range("F7").Value = 'fi lled'
$ws.Range("G7").Value = $false

# Row 8 (id=7)
$ws.Range("B8").Value = 'What was the first time the crew members had been able to find the planet?'
$ws.Range("C8").Value = 'Despite the chaos, there was a certain energy to the city,  a rhythm that kept everything moving, even as night began to fall and the neon lights of the city  flickered to life. As the spaceship approached the unkn own planet, the crew members prepared for their first  descent. They had spent months traveling through the vast reaches of space, charting new star  systems, collecting data, and making discoveries that had never been seen before. But this planet  was differ ent.'
$ws.Range("D8").Value = 'unkn own planet'
$ws.Range("E8").Value = 'a saturn'
$ws.Range("F8").Value = 'unkn own planet'
$ws.Range("G8").Value = $false

# Row 9 (id=8)
$ws.Range("B9").Value = 'What was the atmosphere like?'
$ws.Range("C9").Value = 'It wasn’t just another barren world or distant moon. There was something about it  that felt unique, almost alive. Its atmosphere was thick with an unusual energy, and the surface  appeared to be teeming with strange, uncharted life forms. The crew, cla d in their space suits, felt  a mix of excitement and apprehension as they prepared for their descent into the unknown.'
$ws.Range("D9").Value = 'thick with an unusual energy'
$ws.Range("E9").Value = 'a sweltering night'
$ws.Range("F9").Value = 'thick with an unusual energy'
$ws.Range("G9").Value = $true

# Row 10 (id=9)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'What was the planet''s surface covered in?'
$ws.Range("C10").Value = 'The  planet''s surface was covered in dense forests, vast deserts, and deep, uncharted oceans, each one  holding its own mysteries. As the s pacecraft touched down on the soft soil, the crew knew that  this was only the beginning of a new chapter in their exploration of the universe. What they  would find here could change everything they knew about life beyond Earth, and the adventure  ahead woul d be like none they had experienced before.'
$ws.Range("D10").Value = 'dense forests'
$ws.Range("E10").Value = 'ice'
$ws.Range("F10").Value = 'dense forests'
$ws.Range("G10").Value = $false
